$d = $word.ActiveDocument

# 1) "Sicherung von Benutzereinstellungen/Orten" -> "Sicherung von Orten"
$d.Content.Find.Execute("Sicherung von Benutzereinstellungen/Orten", $true, $false, $false, $false, $false, $true, 1, $false, "Sicherung von Orten", 2)

# 2) Collapse "Nutzung von " + "Github" + " zur gemeinsamen Kollaboration mit dem Team" into one run/text
$d.Content.Find.Execute("Nutzung von Github zur gemeinsamen Kollaboration mit dem Team", $true, $false, $false, $false, $false, $true, 1, $false, "Nutzung von Github zur gemeinsamen Kollaboration mit dem Team", 2)

# 3) Collapse ", " + "Eslem" + " " + "Özdal" into one run/text
$d.Content.Find.Execute(", Eslem Özdal", $true, $false, $false, $false, $false, $true, 1, $false, ", Eslem Özdal", 2)
